# Apply FHIR IG terminology/profile correction edits to the workbook.
#
# Changes (per the regenerated CodeSystem report):
#   - Metadata sheet, "Experimental" row (row 7): the previously-blank
#     Value cell (B7) now carries the literal text "false".
#   - Metadata sheet, "Date" row (row 8): the Value cell (B8) is updated
#     from the old generation timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- B7: "Experimental" -> "false" -------------------------------------
# A bare Value assignment of "false"/"true" is auto-typed as a Boolean by
# Excel (typing FALSE into a cell yields a logical value, not text) -
# same as this report needs, which stores it as literal text. Forcing it
# with a leading apostrophe keeps it text, but also stamps the cell with
# a "quote prefix" format flag/new style. Re-apply the original cell's
# plain formatting (copy format only, from the untouched sibling cell in
# column A on the same row) to drop that quote-prefix style while
# keeping the cell's string value and its original look intact.
$ws.Range("B7").Value = "'false"
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)  # xlPasteFormats

# --- B8: Date timestamp refresh -----------------------------------------
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
